# Update countries & provincias Spain
# Updates the "Datos actualizados" timestamp in A1 and refreshes the
# per-country statistics (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes) for the rows whose
# figures changed in the latest data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "last refreshed" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 12:52"

# Each entry: row number followed by the new values for columns
# B (Casos totales), C (Nuevos casos), D (Casos activos), E (Recuperados),
# F (Casos criticos), G (Muertes hoy), H (Muertes)
$updates = @(
    @(15, 25503, 88, 12700, 11686, 386, 11, 1117),
    @(46, 3231, 252, 334, 2890, 37, 0, 7),
    @(47, 3102, 325, 97, 2912, 45, 10, 93),
    @(48, 3064, 90, 300, 2708, 80, 0, 56),
    @(66, 1348, 212, 591, 751, 4, 0, 6),
    @(67, 1332, 23, 102, 1203, 9, 2, 27),
    @(68, 1300, 66, 150, 1148, 26, 1, 2),
    @(69, 1212, 7, 152, 1005, 34, 2, 55),
    @(103, 384, 6, 44, 337, 4, 0, 3),
    @(108, 291, 11, 178, 111, 1, 0, 2),
    @(109, 290, 0, 58, 230, 0, 0, 2),
    @(115, 232, 4, 131, 99, 13, 0, 2),
    @(117, 214, 4, 56, 151, 1, 0, 7),
    @(127, 136, 0, 107, 28, 2, 0, 1)
)

foreach ($entry in $updates) {
    $row = $entry[0]
    $col = 2
    for ($i = 1; $i -lt $entry.Length; $i++) {
        $ws.Cells.Item($row, $col).Value = $entry[$i]
        $col = $col + 1
    }
}
